$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("Damian Lillard", "PG", "Milwaukee Bucks")
    ,@("Gary Trent Jr.", "PG,SG,SF", "Milwaukee Bucks")
    ,@("Derrick White", "PG,SG", "Boston Celtics")
    ,@("Cade Cunningham", "PG,SG", "Detroit Pistons")
    ,@("Ausar Thompson", "SF,PF", "Detroit Pistons")
    ,@("Naz Reid", "PF,C", "Minnesota Timberwolves")
    ,@("Shaedon Sharpe", "SG,SF", "Portland Trail Blazers")
    ,@("Malik Monk", "PG,SG,SF", "Sacramento Kings")
    ,@("Isaiah Hartenstein", "C", "Oklahoma City Thunder")
    ,@("Zach Edey", "C", "Memphis Grizzlies")
    ,@("LaMelo Ball", "PG,SG", "Charlotte Hornets")
    ,@("Coby White", "PG,SG", "Chicago Bulls")
    ,@("Dalton Knecht", "SG,SF", "Los Angeles Lakers")
    ,@("Aaron Nesmith", "SF,PF", "Indiana Pacers")
    ,@("Onyeka Okongwu", "PF,C", "Atlanta Hawks")
    ,@("Andrew Wiggins", "SF,PF", "Miami Heat")
    ,@("Anthony Davis", "PF,C", "Dallas Mavericks")
    ,@("Collin Sexton", "PG,SG", "Utah Jazz")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
